$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage so that numeric-looking
# strings (e.g. "12.51", "4.85") are not silently converted to numbers by Excel,
# and restore the cell's original style afterwards so no formatting changes leak in.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "35.624.27"
$ws.Range("E2").Value = "  +2.04%  "
Set-TextValue "D3" "1.908.37"
$ws.Range("E3").Value = "  +3.61%  "
$ws.Range("E4").Value = "  +0.57%  "
Set-TextValue "D5" "245.40"
$ws.Range("E5").Value = "  +5.73%  "
Set-TextValue "D6" "0.635"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("E7").Value = "  +0.37%  "
Set-TextValue "D8" "42.60"
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("E9").Value = "  +3.79%  "
Set-TextValue "D10" "0.0708"
$ws.Range("E10").Value = "  +2.58%  "
Set-TextValue "D11" "0.0996"
$ws.Range("E11").Value = "  +1.43%  "
Set-TextValue "D12" "2.183.82"
$ws.Range("E12").Value = "  +3.52%  "
Set-TextValue "D13" "12.51"
$ws.Range("E13").Value = "  +10.09%  "
Set-TextValue "B14" "Polygon"
Set-TextValue "C14" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D14" "0.694"
$ws.Range("E14").Value = "  +3.54%  "
Set-TextValue "B15" "WrappedEther"
Set-TextValue "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D15" "1.886.71"
$ws.Range("E15").Value = "  +2.40%  "
Set-TextValue "D16" "4.85"
$ws.Range("E16").Value = "  +4.17%  "
Set-TextValue "D17" "35.610.01"
$ws.Range("E17").Value = "  +2.01%  "
Set-TextValue "D18" "72.11"
$ws.Range("E18").Value = "  +3.24%  "
$ws.Range("E19").Value = "  +2.75%  "
Set-TextValue "D20" "244.70"
$ws.Range("E20").Value = "  +1.94%  "
Set-TextValue "D21" "12.50"
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("E22").Value = "  +4.01%  "
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  +1.65%  "
Set-TextValue "D25" "171.40"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  +31.46%  "
Set-TextValue "D27" "8.52"
$ws.Range("E27").Value = "  +8.36%  "
Set-TextValue "D28" "18.00"
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("E29").Value = "  +1.98%  "
Set-TextValue "D30" "0.956"
$ws.Range("E30").Value = "  +28.90%  "
Set-TextValue "D31" "4.11"
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("E33").Value = "  +0.57%  "
Set-TextValue "D34" "4.13"
$ws.Range("E34").Value = "  +6.09%  "
Set-TextValue "D35" "1.75"
$ws.Range("E35").Value = "  +7.54%  "
Set-TextValue "D36" "2.05"
$ws.Range("E36").Value = "  +4.97%  "
$ws.Range("E37").Value = "  +6.45%  "
$ws.Range("E38").Value = "  +4.53%  "
$ws.Range("E39").Value = "  +4.66%  "
Set-TextValue "D40" "91.65"
$ws.Range("E40").Value = "  +2.09%  "
Set-TextValue "D41" "1.361.17"
$ws.Range("E41").Value = "  +1.54%  "
Set-TextValue "D42" "15.25"
$ws.Range("E42").Value = "  +4.83%  "
Set-TextValue "D43" "0.0598"
$ws.Range("E43").Value = "  +12.96%  "
$ws.Range("E44").Value = "  +4.92%  "
Set-TextValue "D45" "12.91"
$ws.Range("E45").Value = "  +14.20%  "
Set-TextValue "D46" "47.30"
$ws.Range("E46").Value = "  +39.85%  "
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("E49").Value = "  +5.82%  "
Set-TextValue "D50" "2.095.66"
$ws.Range("E50").Value = "  +3.45%  "
Set-TextValue "D51" "3.53"
$ws.Range("E51").Value = "  +4.07%  "
